# Generate Report for Handback
#
# The file "e22cf724-ad23-4eda-9f93-9afdb7c690f3.md" has been handed back
# (in sync with en-US) for both the zh-cn and de-de locales. Update the
# Overview sheet's status for that row, update the per-locale detail
# sheets' status, and record the new "Latest Handback DateTime" for each
# locale.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# Overview sheet: row 3 corresponds to e22cf724-ad23-4eda-9f93-9afdb7c690f3.md
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# zh-cn detail sheet: row 3 is the same file; update status and handback time
$zhcn.Range("B3").Value = $handedBack
$zhcn.Range("G3").Value = "2016-01-27 07:56:31"

# de-de detail sheet: row 3 is the same file; update status and handback time
$dede.Range("B3").Value = $handedBack
$dede.Range("G3").Value = "2016-01-27 07:56:54"
